$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# "Loan RBI, Variable Instalments": split the old "Late" / "Over Due" block by
# inserting a new (blank) column before column N. This shifts the old column N
# ("Late") to O and the old column P ("Over Due") to Q, matching the new
# variable-instalment repayment schedule layout.
$ws.Columns("N").Insert()

# Keep the new column's width in line with its neighbour (column M) instead of
# falling back to the sheet default width.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# "Repayment Schedule" becomes the active/visible tab (previously "Transactions").
$ws.Activate()
[void]$ws.Range("Q6").Select()
